# "last commit of the day" - add a new BOM line (row 12) for a NAND gate
# part, mirroring the pattern of the existing component rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Best.-Nr. / Bezeichnung / Details / Lieferant for the new component
# (written in this order so new shared-string entries land in the same
# slots as the source workbook: 863-... , NAND Gatter , MC74VHC1G135DBVT1G ).
$ws.Range("F12").Value = "863-MC74VHC1G135DBVT "
$ws.Range("B12").Value = "NAND Gatter "
$ws.Range("C12").Value = "MC74VHC1G135DBVT1G "
$ws.Range("E12").Value = "Mouser"

# The new, longer text wraps onto a second line, so the row grows from
# its single-line height (16.5) to a two-line height (33), same as the
# other wrapped rows (e.g. row 4).
$ws.Rows.Item(12).RowHeight = 33

# Last thing the author did before leaving: click into C12.
$ws.Range("C12").Select() | Out-Null
